# Apply "Add data for 2022-06-22" update to the carjacking-by-month-yoy workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-06-14"

# 2. Update the June row label text (A7) to reflect the new "through" date.
$ws.Range("A7").Value = "June (through 06-14)"

# 3. Update the June row (row 7) counts for each year column (B..I).
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 18
$ws.Range("D7").Value = 23
$ws.Range("E7").Value = 27
$ws.Range("F7").Value = 20
$ws.Range("G7").Value = 55
$ws.Range("H7").Value = 51
$ws.Range("I7").Value = 65

# 4. Update the Total row (row 8) counts for each year column (B..I).
$ws.Range("B8").Value = 115
$ws.Range("C8").Value = 227
$ws.Range("D8").Value = 339
$ws.Range("E8").Value = 322
$ws.Range("F8").Value = 224
$ws.Range("G8").Value = 413
$ws.Range("H8").Value = 682
$ws.Range("I8").Value = 728
